$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.577.50"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "2.939.86"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'377.13"
$ws.Range("E5").Value = "  +6.25%  "

$ws.Range("D6").Value = "'104.39"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("E7").Value = "  -2.34%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  -2.87%  "

$ws.Range("D10").Value = "'37.05"
$ws.Range("E10").Value = "  -2.32%  "

$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("D13").Value = "'18.45"
$ws.Range("E13").Value = "  -2.86%  "

$ws.Range("D14").Value = "3.398.38"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("E15").Value = "  -2.25%  "

$ws.Range("D16").Value = "2.927.96"
$ws.Range("E16").Value = "  -1.57%  "

$ws.Range("E17").Value = "  -5.56%  "

$ws.Range("D18").Value = "51.526.42"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").Value = "'3.42"
$ws.Range("E19").Value = "  +1.56%  "

$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").Value = "'13.06"
$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("E23").Value = "  -0.85%  "

$ws.Range("D24").Value = "'262.51"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").Value = "'2.83"
$ws.Range("E25").Value = "  +4.25%  "

$ws.Range("D26").Value = "'4.14"
$ws.Range("E26").Value = "  -4.59%  "

$ws.Range("D27").Value = "'7.18"
$ws.Range("E27").Value = "  +13.45%  "

$ws.Range("E28").Value = "  -5.84%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "'25.89"
$ws.Range("E30").Value = "  -3.61%  "

$ws.Range("D31").Value = "'7.33"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").Value = "'0.103"
$ws.Range("E32").Value = "  -5.90%  "

$ws.Range("D33").Value = "'9.85"
$ws.Range("E33").Value = "  -2.76%  "

$ws.Range("D34").Value = "'51.92"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("E35").Value = "  -3.44%  "

$ws.Range("D36").Value = "'34.23"
$ws.Range("E36").Value = "  -4.93%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("E39").Value = "  -7.10%  "

$ws.Range("D40").Value = "'17.09"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("E41").Value = "  -8.22%  "

$ws.Range("E42").Value = "  -5.28%  "

$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").Value = "'124.76"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("D45").Value = "'21.87"
$ws.Range("E45").Value = "  -5.46%  "

$ws.Range("E46").Value = "  -4.69%  "

$ws.Range("D47").Value = "'0.276"
$ws.Range("E47").Value = "  +14.03%  "

$ws.Range("D48").Value = "2.025.40"
$ws.Range("E48").Value = "  -4.50%  "

$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").Value = "'3.18"
$ws.Range("E50").Value = "  -3.60%  "

$ws.Range("E51").Value = "  -1.76%  "
